$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 123
$ws.Range("B3").Value = "Maria Clara"
$ws.Range("C3").Value = $false

$ws.Range("A4").Value = 111
$ws.Range("B4").Value = "Cecilia Santos"
$ws.Range("C4").Value = $false

$ws.Range("A5").Value = 12
$ws.Range("B5").Value = "Daniel Valente"
$ws.Range("C5").Value = $true
